# Add a new "Flip" sheet (a copy of "Clockwise") with a new, randomized
# ordering of the level numbers, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# The workbook's internal sheetId counter only ever increases (max existing
# sheetId + 1), it is never reused just because a sheet was removed. The
# target file expects the new "Flip" sheet to carry sheetId 6, while a
# straightforward Add()/Copy() right now would hand out sheetId 5 (the
# highest current id, 4, plus one). Burn id 5 on a throwaway sheet first so
# the real "Flip" sheet is minted with id 6, then discard the throwaway.
$throwaway = $wb.Worksheets.Add()

# Clone the "Clockwise" sheet (keeps formatting/column widths/styles
# identical) and drop the clone at the end of the tab strip.
$source = $wb.Worksheets.Item("Clockwise")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)

# Remove the throwaway sheet now that the real copy already grabbed id 6.
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null

# The copy is named "Clockwise (2)" - rename it and make it the active tab.
$flip = $wb.Worksheets.Item("Clockwise (2)")
$flip.Name = "Flip"
$flip.Activate()

# Randomized level order for the new "Flip" sheet.
$order = @(4, 8, 12, 1, 5, 9, 13, 2, 6, 10, 14, 3, 7, 11)
for ($row = 1; $row -le $order.Length; $row++) {
    $flip.Cells.Item($row, 2).Value = $order[$row - 1]
}

# Match the author's last selection on the new sheet.
$flip.Range("C15").Select() | Out-Null
